$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 975, pushing existing rows (975 onward) down by 2.
$ws.Rows.Item(975).Insert()
$ws.Rows.Item(975).Insert()

# --- New row 975 (Primera) ---
$ws.Cells.Item(975, 1).Value = 3
$ws.Cells.Item(975, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(975, 3).Value = "Coquimbo"
$ws.Cells.Item(975, 4).Value = 45021
$ws.Cells.Item(975, 5).Value = 5
$ws.Cells.Item(975, 6).Value = 100112023
$ws.Cells.Item(975, 7).Value = "Brócoli"
$ws.Cells.Item(975, 8).Value = "Sin especificar"
$ws.Cells.Item(975, 9).Value = "Primera"
$ws.Cells.Item(975, 10).Value = 2500
$ws.Cells.Item(975, 11).Value = 900
$ws.Cells.Item(975, 12).Value = 1000
$ws.Cells.Item(975, 13).Value = 952
$ws.Cells.Item(975, 14).Value = "$/unidad"
$ws.Cells.Item(975, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(975, 16).Value = 952
$ws.Cells.Item(975, 17).Value = 1
$ws.Cells.Item(975, 18).Value = "Hortaliza"

# --- New row 976 (Segunda) ---
$ws.Cells.Item(976, 1).Value = 3
$ws.Cells.Item(976, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(976, 3).Value = "Coquimbo"
$ws.Cells.Item(976, 4).Value = 45021
$ws.Cells.Item(976, 5).Value = 5
$ws.Cells.Item(976, 6).Value = 100112023
$ws.Cells.Item(976, 7).Value = "Brócoli"
$ws.Cells.Item(976, 8).Value = "Sin especificar"
$ws.Cells.Item(976, 9).Value = "Segunda"
$ws.Cells.Item(976, 10).Value = 1100
$ws.Cells.Item(976, 11).Value = 800
$ws.Cells.Item(976, 12).Value = 800
$ws.Cells.Item(976, 13).Value = 800
$ws.Cells.Item(976, 14).Value = "$/unidad"
$ws.Cells.Item(976, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(976, 16).Value = 800
$ws.Cells.Item(976, 17).Value = 1
$ws.Cells.Item(976, 18).Value = "Hortaliza"
